# Update Shadow Rate with Latest Data
# Refresh the fedfundsrate_shadow series (column C) with newly recomputed
# shadow-rate estimates, and correct the 1984-Q1 observation (row 2) whose
# upstream source value changed slightly on the refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (1984): both fedfundsrate and fedfundsrate_shadow refreshed together
$ws.Range("B2").Value = 9.6866666666667545
$ws.Range("C2").Value = 9.6866666666667545

# 2009-2012 ZIRP period: updated shadow-rate estimates
$ws.Range("C102").Value = 1.5775537756388847
$ws.Range("C103").Value = 0.16775172713514586
$ws.Range("C104").Value = -0.74389033782763514
$ws.Range("C105").Value = -0.80245942206732535
$ws.Range("C106").Value = -0.63026216887094311
$ws.Range("C107").Value = -2.1024210675325117
$ws.Range("C108").Value = -1.9322171886070216
$ws.Range("C109").Value = -2.5774161397115525
$ws.Range("C110").Value = -2.0215396513580952
$ws.Range("C111").Value = -1.6219739495136776
$ws.Range("C112").Value = -2.9124924104203531
$ws.Range("C113").Value = -2.479521876826718
$ws.Range("C114").Value = -3.3719517909621466
$ws.Range("C115").Value = -3.1497286216460307
$ws.Range("C116").Value = -2.7029206466402056
$ws.Range("C117").Value = -3.9116689347991951
$ws.Range("C118").Value = -2.3761542278954995
$ws.Range("C119").Value = -1.6485086014214279
$ws.Range("C120").Value = -1.2481728989292451
$ws.Range("C121").Value = -1.3313759744145282
$ws.Range("C122").Value = -1.3766592037401049
$ws.Range("C123").Value = -1.3059446539577801
$ws.Range("C124").Value = -0.89056223411516644
$ws.Range("C125").Value = -0.41933499758491921
$ws.Range("C126").Value = 0.12215563006285901
$ws.Range("C127").Value = 0.012575595825126129
$ws.Range("C128").Value = 0.044410665622063306

# 2020-2022 ZIRP period: updated shadow-rate estimates (latest data pull)
$ws.Range("C147").Value = 0.11290968855326433
$ws.Range("C148").Value = -11.399859306774584
$ws.Range("C149").Value = -4.9636593258756356
$ws.Range("C150").Value = -4.0560136448109407
$ws.Range("C151").Value = -3.887693620518673
$ws.Range("C152").Value = -2.6368599388076897
$ws.Range("C153").Value = -1.4418539797378283
$ws.Range("C154").Value = -0.47727702503747027
